# Updated symbol list (cryptocurrency price/volume snapshot refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are stored as literal text (e.g. "244.58", with
# significant trailing zeros such as "0.00002100"), so force text formatting
# before assigning, otherwise Excel would silently coerce the numeric-looking
# strings into floating point numbers and lose the exact textual form.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.58"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.86"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06041"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.395"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8135"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9223"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1438"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07494"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03383"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03048"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09401"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001601"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04796"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0005940"
$ws.Range("E17").Value = "16OneONE"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.005408"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.004160"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0009889"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.00008803"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.651"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.439"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1304"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002900"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03988"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.003036"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002721"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.006376"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005255"

$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002526"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
